$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.334.67'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '2.092.01'
$ws.Range("E3").Value = '  +4.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5228'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4343'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08816'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.164'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.94%  '

$ws.Range("D13").Value = '2.089.63'
$ws.Range("E13").Value = '  +3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.736'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.751'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.94%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001128'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06639'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.93'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.340'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.13%  '

$ws.Range("D23").Value = '30.393.98'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("E24").Value = '  +4.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.307'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.62%  '

$ws.Range("D26").Value = '2.331.07'
$ws.Range("E26").Value = '  +3.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.596'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.213'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1071'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.675'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +23.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.177'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.893'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.55%  '

$ws.Range("E36").Value = '  +9.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02586'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06699'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.474'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2262'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6825'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.249'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.50%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6372'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.209'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.615'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.249'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.192'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.40%  '
